# Insert a new weekly data row for "Arándano (blue)" at Vega Monumental
# Concepción. The new row is inserted at row 35 (pushing the existing
# rows 35-81 down to 36-82), matching the target diff where the sheet's
# used range grows from A1:T81 to A1:T82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("35").Insert()

$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value = "2022-01-18"
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100101
$ws.Range("H35").Value = "Berries"
$ws.Range("I35").Value = 100101001
$ws.Range("J35").Value = "Arándano (blue)"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 250
$ws.Range("N35").Value = 3500
$ws.Range("O35").Value = 4000
$ws.Range("P35").Value = 3700
$ws.Range("Q35").Value = "$/bandeja 2 kilos"
$ws.Range("R35").Value = "Provincia de Curicó"
$ws.Range("S35").Value = 1850
$ws.Range("T35").Value = 2
